{"js": "// Add the repo link into the \"Status\" cell of the \"Source code uploaded to\n// GitHub\" row in the CHECKLIST table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// Find the checklist table: header row reads SNo / Item / Status.\nlet target = null;\nfor (let i = 0; i < tables.items.length; i++) {\n  const t = tables.items[i];\n  t.rows.load(\"items\");\n  await context.sync();\n  if (t.rows.items.length === 0) continue;\n  const headerRow = t.rows.items[0];\n  headerRow.cells.load(\"items\");\n  await context.sync();\n  if (headerRow.cells.items.length < 3) continue;\n  headerRow.cells.items[0].body.load(\"text\");\n  headerRow.cells.items[1].body.load(\"text\");\n  headerRow.cells.items[2].body.load(\"text\");\n  await context.sync();\n  const h0 = headerRow.cells.items[0].body.text.trim();\n  const h1 = headerRow.cells.items[1].body.text.trim();\n  const h2 = headerRow.cells.items[2].body.text.trim();\n  if (h0 === \"SNo\" && h1 === \"Item\" && h2 === \"Status\") {\n    target = t;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the CHECKLIST table.\");\n}\n\ntarget.rows.load(\"items\");\nawait context.sync();\n\n// Find the row whose \"Item\" cell mentions uploading the source code to GitHub.\nlet statusCell = null;\nfor (let r = 0; r < target.rows.items.length; r++) {\n  const row = target.rows.items[r];\n  row.cells.load(\"items\");\n  await context.sync();\n  if (row.cells.items.length < 3) continue;\n  const itemCell = row.cells.items[1];\n  itemCell.body.load(\"text\");\n  await context.sync();\n  if (itemCell.body.text.indexOf(\"Source code uploaded to\") !== -1) {\n    statusCell = row.cells.items[2];\n    break;\n  }\n}\n\nif (!statusCell) {\n  throw new Error(\"Could not locate the GitHub upload row.\");\n}\n\nstatusCell.body.insertText(\"https://github.com/Sahith02/YACS\", \"Replace\");\nawait context.sync();\n\n// Word auto-fits this table's column grid to its (now longer) content;\n// mirror the resulting column widths (twips -> points, 20 twips per pt).\ntarget.rows.load(\"items\");\nawait context.sync();\nconst newWidthsTwips = [1059, 3982, 3589];\nfor (let r = 0; r < target.rows.items.length; r++) {\n  const row = target.rows.items[r];\n  row.cells.load(\"items\");\n  await context.sync();\n  for (let c = 0; c < row.cells.items.length && c < newWidthsTwips.length; c++) {\n    row.cells.items[c].columnWidth = newWidthsTwips[c] / 20;\n  }\n}\nawait context.sync();\n", "ps1": "# Add the repo link into the \"Status\" cell of the \"Source code uploaded to\n# GitHub\" row in the CHECKLIST table.\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    if ($t.Columns.Count -lt 3) { continue }\n    $h0 = ($t.Cell(1, 1).Range.Text -replace \"[\\r\\a]\", \"\").Trim()\n    $h1 = ($t.Cell(1, 2).Range.Text -replace \"[\\r\\a]\", \"\").Trim()\n    $h2 = ($t.Cell(1, 3).Range.Text -replace \"[\\r\\a]\", \"\").Trim()\n    if ($h0 -eq \"SNo\" -and $h1 -eq \"Item\" -and $h2 -eq \"Status\") {\n        $target = $t\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not locate the CHECKLIST table.\"\n}\n\n$targetRow = 0\nfor ($r = 1; $r -le $target.Rows.Count; $r++) {\n    $itemText = $target.Cell($r, 2).Range.Text\n    if ($itemText -like \"*Source code uploaded to*\") {\n        $targetRow = $r\n        break\n    }\n}\n\nif ($targetRow -eq 0) {\n    throw \"Could not locate the GitHub upload row.\"\n}\n\n$target.Cell($targetRow, 3).Range.Text = \"https://github.com/Sahith02/YACS\"\n\n# Word auto-fits this table's column grid to its (now longer) content;\n# mirror the resulting column widths (twips -> points, 20 twips per pt).\n$newWidthsTwips = @(1059, 3982, 3589)\nfor ($r = 1; $r -le $target.Rows.Count; $r++) {\n    for ($c = 1; $c -le $target.Columns.Count; $c++) {\n        $target.Cell($r, $c).Width = $newWidthsTwips[$c - 1] / 20\n    }\n}\n"}
